$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 4 ("Modesto"/"Primera" 44938 entry),
# shifting the existing rows 4-17 down to 6-19.
$ws.Range("A4:A5").EntireRow.Insert()

# New row 4: Modesto / Primera, 16 kilos empedrada box, week of 44944
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44944
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103003
$ws.Range("J4").Value = "Damasco"
$ws.Range("K4").Value = "Modesto"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("Q4").Value = "$/caja 16 kilos empedrada"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1031
$ws.Range("T4").Value = 16

# New row 5: Modesto / Segunda, 16 kilos empedrada box, week of 44944
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44944
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100103
$ws.Range("H5").Value = "Frutos de hueso (carozo)"
$ws.Range("I5").Value = 100103003
$ws.Range("J5").Value = "Damasco"
$ws.Range("K5").Value = "Modesto"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 14000
$ws.Range("Q5").Value = "$/caja 16 kilos empedrada"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 875
$ws.Range("T5").Value = 16
